$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename Id columns (drop leading '#')
$ws.Range("A2").Value = "Id"
$ws.Range("A3").Value = "PlayerId"

# Clear the ProtocolType ("E") column data below the header, and the
# leftover "fk" marker in F4 (StructureItemId row), since the Packet
# protocol annotations are no longer part of this model.
$ws.Range("E2:E12").ClearContents()
$ws.Range("F4").ClearContents()

# Reset the active selection as left by the author after editing.
$ws.Range("I6").Select()
